$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: add idx value in A4
$ws.Range("A4").Value = 2

# Row 5: add idx value in A5
$ws.Range("A5").Value = 3

# New row 6: sphere/anchor4 color data (idx=4, tag coords, anchors, dists)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 5

$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Formula = "=SQRT((E6-`$B6)^2+(F6-`$C6)^2+(G6-`$D6)^2)"

$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 10
$ws.Range("K6").Formula = "=10/3"
$ws.Range("L6").Formula = "=SQRT((I6-`$B6)^2+(J6-`$C6)^2+(K6-`$D6)^2)"

$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 10
$ws.Range("O6").Formula = "=2*(10/3)"
$ws.Range("P6").Formula = "=SQRT((M6-`$B6)^2+(N6-`$C6)^2+(O6-`$D6)^2)"

$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 10
$ws.Range("T6").Formula = "=SQRT((Q6-`$B6)^2+(R6-`$C6)^2+(S6-`$D6)^2)"

# Apply same style (s=7) as other columns A-D cells by copying format from A3
$ws.Range("A3").Copy()
$ws.Range("A4:A6").PasteSpecial(-4122) # xlPasteFormats

# Copy full row style from row 5 into row 6, then re-apply values (since paste may overwrite values)
$ws.Range("B5:T5").Copy()
$ws.Range("B6:T6").PasteSpecial(-4122) # xlPasteFormats

# Selection / view changes
$ws.Range("Q7").Select()

$wb.Save()
